$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "C3"
$ws.Cells.Item(2,3).Value = "C3ar1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 34.739995
$ws.Cells.Item(2,8).Value = 104.219985
$ws.Cells.Item(2,9).Value = 0.1827267341390226
$ws.Cells.Item(2,10).Value = 0.1827267341390226
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 11.838451
$ws.Cells.Item(2,14).Value = 35.515353
$ws.Cells.Item(2,15).Value = 0.9341683151774448
$ws.Cells.Item(2,16).Value = 0.9341683151774446
$ws.Cells.Item(2,17).Value = 411.267728547745
$ws.Cells.Item(2,18).Value = 3701.409556929705
$ws.Cells.Item(2,19).Value = 0.1706975253685276
$ws.Cells.Item(2,20).Value = 0.1706975253685276

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "C3"
$ws.Cells.Item(3,3).Value = "C3ar1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 34.739995
$ws.Cells.Item(3,8).Value = 104.219985
$ws.Cells.Item(3,9).Value = 0.1827267341390226
$ws.Cells.Item(3,10).Value = 0.1827267341390226
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.623012
$ws.Cells.Item(3,14).Value = 1.869036
$ws.Cells.Item(3,15).Value = 0.0491616741392375
$ws.Cells.Item(3,16).Value = 0.04916167413923749
$ws.Cells.Item(3,17).Value = 21.64343376494
$ws.Cells.Item(3,18).Value = 194.79090388446
$ws.Cells.Item(3,19).Value = 0.008983152160269715
$ws.Cells.Item(3,20).Value = 0.008983152160269712

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "C3"
$ws.Cells.Item(4,3).Value = "C3ar1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 34.739995
$ws.Cells.Item(4,8).Value = 104.219985
$ws.Cells.Item(4,9).Value = 0.1827267341390226
$ws.Cells.Item(4,10).Value = 0.1827267341390226
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.2112543333333333
$ws.Cells.Item(4,14).Value = 0.633763
$ws.Cells.Item(4,15).Value = 0.0166700106833178
$ws.Cells.Item(4,16).Value = 0.0166700106833178
$ws.Cells.Item(4,17).Value = 7.338974483728333
$ws.Cells.Item(4,18).Value = 66.050770353555
$ws.Cells.Item(4,19).Value = 0.003046056610225279
$ws.Cells.Item(4,20).Value = 0.003046056610225279

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "C3"
$ws.Cells.Item(5,3).Value = "C3ar1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 154.8642143333334
$ws.Cells.Item(5,8).Value = 464.5926430000001
$ws.Cells.Item(5,9).Value = 0.8145606273154508
$ws.Cells.Item(5,10).Value = 0.8145606273154508
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 11.838451
$ws.Cells.Item(5,14).Value = 35.515353
$ws.Cells.Item(5,15).Value = 0.9341683151774448
$ws.Cells.Item(5,16).Value = 0.9341683151774446
$ws.Cells.Item(5,17).Value = 1833.352413038665
$ws.Cells.Item(5,18).Value = 16500.17171734798
$ws.Cells.Item(5,19).Value = 0.7609367288291572
$ws.Cells.Item(5,20).Value = 0.760936728829157

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "C3"
$ws.Cells.Item(6,3).Value = "C3ar1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 154.8642143333334
$ws.Cells.Item(6,8).Value = 464.5926430000001
$ws.Cells.Item(6,9).Value = 0.8145606273154508
$ws.Cells.Item(6,10).Value = 0.8145606273154508
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.623012
$ws.Cells.Item(6,14).Value = 1.869036
$ws.Cells.Item(6,15).Value = 0.0491616741392375
$ws.Cells.Item(6,16).Value = 0.04916167413923749
$ws.Cells.Item(6,17).Value = 96.48226390023869
$ws.Cells.Item(6,18).Value = 868.3403751021481
$ws.Cells.Item(6,19).Value = 0.04004516412673507
$ws.Cells.Item(6,20).Value = 0.04004516412673506

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "C3"
$ws.Cells.Item(7,3).Value = "C3ar1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 154.8642143333334
$ws.Cells.Item(7,8).Value = 464.5926430000001
$ws.Cells.Item(7,9).Value = 0.8145606273154508
$ws.Cells.Item(7,10).Value = 0.8145606273154508
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.2112543333333333
$ws.Cells.Item(7,14).Value = 0.633763
$ws.Cells.Item(7,15).Value = 0.0166700106833178
$ws.Cells.Item(7,16).Value = 0.0166700106833178
$ws.Cells.Item(7,17).Value = 32.71573635617878
$ws.Cells.Item(7,18).Value = 294.4416272056091
$ws.Cells.Item(7,19).Value = 0.01357873435955861
$ws.Cells.Item(7,20).Value = 0.01357873435955861

$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "C3"
$ws.Cells.Item(8,3).Value = "C3ar1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.5157266666666667
$ws.Cells.Item(8,8).Value = 1.54718
$ws.Cells.Item(8,9).Value = 0.002712638545526686
$ws.Cells.Item(8,10).Value = 0.002712638545526686
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 11.838451
$ws.Cells.Item(8,14).Value = 35.515353
$ws.Cells.Item(8,15).Value = 0.9341683151774448
$ws.Cells.Item(8,16).Value = 0.9341683151774446
$ws.Cells.Item(8,17).Value = 6.105404872726666
$ws.Cells.Item(8,18).Value = 54.94864385453999
$ws.Cells.Item(8,19).Value = 0.002534060979760059
$ws.Cells.Item(8,20).Value = 0.002534060979760058

$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "C3"
$ws.Cells.Item(9,3).Value = "C3ar1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.5157266666666667
$ws.Cells.Item(9,8).Value = 1.54718
$ws.Cells.Item(9,9).Value = 0.002712638545526686
$ws.Cells.Item(9,10).Value = 0.002712638545526686
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.623012
$ws.Cells.Item(9,14).Value = 1.869036
$ws.Cells.Item(9,15).Value = 0.0491616741392375
$ws.Cells.Item(9,16).Value = 0.04916167413923749
$ws.Cells.Item(9,17).Value = 0.3213039020533333
$ws.Cells.Item(9,18).Value = 2.89173511848
$ws.Cells.Item(9,19).Value = 0.0001333578522327181
$ws.Cells.Item(9,20).Value = 0.0001333578522327181

$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "C3"
$ws.Cells.Item(10,3).Value = "C3ar1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.5157266666666667
$ws.Cells.Item(10,8).Value = 1.54718
$ws.Cells.Item(10,9).Value = 0.002712638545526686
$ws.Cells.Item(10,10).Value = 0.002712638545526686
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.2112543333333333
$ws.Cells.Item(10,14).Value = 0.633763
$ws.Cells.Item(10,15).Value = 0.0166700106833178
$ws.Cells.Item(10,16).Value = 0.0166700106833178
$ws.Cells.Item(10,17).Value = 0.1089494931488889
$ws.Cells.Item(10,18).Value = 0.98054543834
$ws.Cells.Item(10,19).Value = 0.00004521971353390951
$ws.Cells.Item(10,20).Value = 0.00004521971353390952

